$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 17 solution run times: Part 1 (B21) and Part 2 (C21)
$ws.Range("B21").Value = 3.02397709997603
$ws.Range("C21").Value = 0.00061240000650286599

# Move the active selection to the Total cell for day 17
$ws.Range("E21").Select()
